$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.956.29'

$ws.Range('D3').Value = '1.636.25'
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = "'212.30"
$ws.Range('E5').Value = '  -0.79%  '

$ws.Range('E6').Value = '  -1.03%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('E8').Value = '  -1.25%  '

$ws.Range('E9').Value = '  -2.72%  '

$ws.Range('E10').Value = '  +0.03%  '

$ws.Range('E11').Value = '  +1.05%  '

$ws.Range('D12').Value = '1.868.19'
$ws.Range('E12').Value = '  -0.88%  '

$ws.Range('D13').Value = '1.635.19'
$ws.Range('E13').Value = '  -0.93%  '

$ws.Range('E14').Value = '  -0.58%  '

$ws.Range('D15').Value = "'0.566"
$ws.Range('E15').Value = '  -0.24%  '

$ws.Range('D16').Value = "'65.45"
$ws.Range('E16').Value = '  -0.47%  '

$ws.Range('D17').Value = '27.962.90'
$ws.Range('E17').Value = '  -0.30%  '

$ws.Range('D18').Value = "'230.77"
$ws.Range('E18').Value = '  -0.89%  '

$ws.Range('D19').Value = '0.0₃0725'
$ws.Range('E19').Value = '  -0.07%  '

$ws.Range('D20').Value = "'7.56"
$ws.Range('E20').Value = '  -1.52%  '

$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('E22').Value = '  -0.73%  '

$ws.Range('D23').Value = "'10.39"
$ws.Range('E23').Value = '  -3.08%  '

$ws.Range('D24').Value = "'2.07"
$ws.Range('E24').Value = '  -3.93%  '

$ws.Range('D25').Value = "'154.92"
$ws.Range('E25').Value = '  +1.56%  '

$ws.Range('E26').Value = '  +0.63%  '

$ws.Range('E27').Value = '  -0.77%  '

$ws.Range('E28').Value = '  -0.98%  '

$ws.Range('E29').Value = '  -0.12%  '

$ws.Range('E30').Value = '  -0.76%  '

$ws.Range('E31').Value = '  -0.42%  '

$ws.Range('D32').Value = "'3.41"
$ws.Range('E32').Value = '  +1.83%  '

$ws.Range('D33').Value = '1.407.97'
$ws.Range('E33').Value = '  -2.76%  '

$ws.Range('E34').Value = '  -0.21%  '

$ws.Range('E35').Value = '  -0.16%  '

$ws.Range('E36').Value = '  +9.25%  '

$ws.Range('E37').Value = '  +1.36%  '

$ws.Range('E38').Value = '  +0.45%  '

$ws.Range('D39').Value = "'0.563"
$ws.Range('E39').Value = '  +0.86%  '

$ws.Range('E40').Value = '  -2.38%  '

$ws.Range('E41').Value = '  +0.37%  '

$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('D43').Value = "'66.85"
$ws.Range('E43').Value = '  -3.73%  '

$ws.Range('E44').Value = '  +2.13%  '

$ws.Range('D45').Value = "'1.84"
$ws.Range('E45').Value = '  +0.38%  '

$ws.Range('E46').Value = '  -0.98%  '

$ws.Range('E47').Value = '  -0.95%  '

$ws.Range('D48').Value = "'87.99"
$ws.Range('E48').Value = '  -1.37%  '

$ws.Range('E49').Value = '  +1.42%  '

$ws.Range('E50').Value = '  -1.29%  '

$ws.Range('E51').Value = '  -0.40%  '
